# Commit: "added new sample map"
#
# 1. Rename "Sheet0" -> "Sheet00" and replace its grid (was A1:E6) with a new,
#    larger A1:G8 sample map.
# 2. Insert a brand-new worksheet "Sheet01" right after "Sheet00" -- a second
#    (very similar) sample map, also A1:G8. We get it by copying Sheet00
#    (so it inherits the same namespaces / formatting Excel would normally
#    carry over) and then overwriting its values.
# All the other worksheets (Sheet1 .. Sheet12, DEMO) are left untouched;
# Excel renumbers their internal part names (sheetN.xml) automatically once
# the new sheet is spliced in at position 2.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: rename Sheet0 -> Sheet00
# ---------------------------------------------------------------------------
$sheet00 = $wb.Worksheets.Item("Sheet0")
$sheet00.Name = "Sheet00"

# ---------------------------------------------------------------------------
# Step 2: insert the new "Sheet01" worksheet right after Sheet00 by copying
# it (keeps xml namespaces / phoneticPr that a bare .Add() wouldn't have).
# ---------------------------------------------------------------------------
$sheet00.Copy([System.Reflection.Missing]::Value, $sheet00)
$sheet01 = $wb.Worksheets.Item(2)
$sheet01.Name = "Sheet01"

# ---------------------------------------------------------------------------
# Grid data: both maps share the same header and are identical except for
# data row 4 (spreadsheet row 4 = 3rd data row below the header).
# ---------------------------------------------------------------------------
$headerValues = @("col0","col1","col2","col3","col4","col5","col6")
$header = New-Object 'object[,]' 1,7
for ($c = 0; $c -lt 7; $c++) { $header[0,$c] = $headerValues[$c] }

$rows00 = @(
    @(1,1,1,1,1,1,1),
    @(1,-255,0,0,0,0,1),
    @(1,0,0,0,0,0,1),
    @(1,0,0,0,0,0,1),
    @(1,0,0,0,255,0,1),
    @(1,0,0,0,0,0,1),
    @(1,1,1,1,1,1,1)
)
$data00 = New-Object 'object[,]' 7,7
for ($r = 0; $r -lt 7; $r++) {
    for ($c = 0; $c -lt 7; $c++) {
        $data00[$r,$c] = $rows00[$r][$c]
    }
}

$rows01 = @(
    @(1,1,1,1,1,1,1),
    @(1,-255,0,0,0,0,1),
    @(1,1,1,1,1,0,1),
    @(1,0,0,0,0,0,1),
    @(1,0,0,0,255,0,1),
    @(1,0,0,0,0,0,1),
    @(1,1,1,1,1,1,1)
)
$data01 = New-Object 'object[,]' 7,7
for ($r = 0; $r -lt 7; $r++) {
    for ($c = 0; $c -lt 7; $c++) {
        $data01[$r,$c] = $rows01[$r][$c]
    }
}

# ---------------------------------------------------------------------------
# Write Sheet00 (grows from A1:E6 to A1:G8)
# ---------------------------------------------------------------------------
$sheet00.Range("A1:G1").Value = $header
$sheet00.Range("A2:G8").Value = $data00

# ---------------------------------------------------------------------------
# Write Sheet01
# ---------------------------------------------------------------------------
$sheet01.Range("A1:G1").Value = $header
$sheet01.Range("A2:G8").Value = $data01
$sheet01.Range("A1:G8").Select() | Out-Null

# Re-activate Sheet00 so it is the selected/visible tab, matching the
# original workbook's behaviour of showing Sheet0 first.
$sheet00.Activate()
$sheet00.Range("A1").Select() | Out-Null
